$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:E2")
$rng.NumberFormat = "@"
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '29.901.32'
$ws.Range("E2").Value = '  -0.54%  '

$rng = $ws.Range("B3:E3")
$rng.NumberFormat = "@"
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '1.872.39'
$ws.Range("E3").Value = '  -1.04%  '

$rng = $ws.Range("B4:E4")
$rng.NumberFormat = "@"
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.01%  '

$rng = $ws.Range("B5:E5")
$rng.NumberFormat = "@"
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '240.47'
$ws.Range("E5").Value = '  -3.46%  '

$rng = $ws.Range("B6:E6")
$rng.NumberFormat = "@"
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.08%  '

$rng = $ws.Range("B7:E7")
$rng.NumberFormat = "@"
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.4945'
$ws.Range("E7").Value = '  -1.20%  '

$rng = $ws.Range("B8:E8")
$rng.NumberFormat = "@"
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").Value = '44.57'
$ws.Range("E8").Value = '  -2.72%  '

$rng = $ws.Range("B9:E9")
$rng.NumberFormat = "@"
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.2900'
$ws.Range("E9").Value = '  +1.11%  '

$rng = $ws.Range("B10:E10")
$rng.NumberFormat = "@"
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.06557'
$ws.Range("E10").Value = '  -0.12%  '

$rng = $ws.Range("B11:E11")
$rng.NumberFormat = "@"
$ws.Range("B11").Value = 'WrappedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D11").Value = '1.871.60'
$ws.Range("E11").Value = '  -1.07%  '

$rng = $ws.Range("B12:E12")
$rng.NumberFormat = "@"
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '16.62'
$ws.Range("E12").Value = '  -3.58%  '

$rng = $ws.Range("B13:E13")
$rng.NumberFormat = "@"
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07173'
$ws.Range("E13").Value = '  -0.79%  '

$rng = $ws.Range("B14:E14")
$rng.NumberFormat = "@"
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '86.07'
$ws.Range("E14").Value = '  +1.09%  '

$rng = $ws.Range("B15:E15")
$rng.NumberFormat = "@"
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.6560'
$ws.Range("E15").Value = '  -1.96%  '

$rng = $ws.Range("B16:E16")
$rng.NumberFormat = "@"
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = '4.826'
$ws.Range("E16").Value = '  -0.08%  '

$rng = $ws.Range("B17:E17")
$rng.NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '29.897.79'
$ws.Range("E17").Value = '  -0.60%  '

$rng = $ws.Range("B18:E18")
$rng.NumberFormat = "@"
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.000007794'
$ws.Range("E18").Value = '  +3.29%  '

$rng = $ws.Range("B19:E19")
$rng.NumberFormat = "@"
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '0.9993'
$ws.Range("E19").Value = '  -0.07%  '

$rng = $ws.Range("B20:E20")
$rng.NumberFormat = "@"
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '12.65'
$ws.Range("E20").Value = '  -1.97%  '

$rng = $ws.Range("B21:E21")
$rng.NumberFormat = "@"
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.113.41'
$ws.Range("E21").Value = '  -0.69%  '

$rng = $ws.Range("B22:E22")
$rng.NumberFormat = "@"
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.01%  '

$rng = $ws.Range("B23:E23")
$rng.NumberFormat = "@"
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '4.731'
$ws.Range("E23").Value = '  -0.97%  '

$rng = $ws.Range("B24:E24")
$rng.NumberFormat = "@"
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '5.565'
$ws.Range("E24").Value = '  +0.16%  '

$rng = $ws.Range("B25:E25")
$rng.NumberFormat = "@"
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '9.011'
$ws.Range("E25").Value = '  -0.45%  '

$rng = $ws.Range("B26:E26")
$rng.NumberFormat = "@"
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '150.08'
$ws.Range("E26").Value = '  +3.39%  '

$rng = $ws.Range("B27:E27")
$rng.NumberFormat = "@"
$ws.Range("B27").Value = 'BitcoinCash'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D27").Value = '133.76'
$ws.Range("E27").Value = '  -1.57%  '

$rng = $ws.Range("B28:E28")
$rng.NumberFormat = "@"
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '16.64'
$ws.Range("E28").Value = '  -0.96%  '

$rng = $ws.Range("B29:E29")
$rng.NumberFormat = "@"
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '1.891'
$ws.Range("E29").Value = '  -3.38%  '

$rng = $ws.Range("B30:E30")
$rng.NumberFormat = "@"
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.372'
$ws.Range("E30").Value = '  +0.02%  '

$rng = $ws.Range("B31:E31")
$rng.NumberFormat = "@"
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '4.134'
$ws.Range("E31").Value = '  -1.61%  '

$rng = $ws.Range("B32:E32")
$rng.NumberFormat = "@"
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '0.08666'
$ws.Range("E32").Value = '  -0.15%  '

$rng = $ws.Range("B33:E33")
$rng.NumberFormat = "@"
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '3.921'
$ws.Range("E33").Value = '  -0.27%  '

$rng = $ws.Range("B34:E34")
$rng.NumberFormat = "@"
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.05019'
$ws.Range("E34").Value = '  -0.56%  '

$rng = $ws.Range("B35:E35")
$rng.NumberFormat = "@"
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.6952'
$ws.Range("E35").Value = '  +0.21%  '

$rng = $ws.Range("B36:E36")
$rng.NumberFormat = "@"
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.088'
$ws.Range("E36").Value = '  -4.62%  '

$rng = $ws.Range("B37:E37")
$rng.NumberFormat = "@"
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.660'
$ws.Range("E37").Value = '  -0.96%  '

$rng = $ws.Range("B38:E38")
$rng.NumberFormat = "@"
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.682'
$ws.Range("E38").Value = '  -3.29%  '

$rng = $ws.Range("B39:E39")
$rng.NumberFormat = "@"
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '2.149'
$ws.Range("E39").Value = '  -5.81%  '

$rng = $ws.Range("B40:E40")
$rng.NumberFormat = "@"
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01688'
$ws.Range("E40").Value = '  +2.72%  '

$rng = $ws.Range("B41:E41")
$rng.NumberFormat = "@"
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '0.9228'
$ws.Range("E41").Value = '  -4.19%  '

$rng = $ws.Range("B42:E42")
$rng.NumberFormat = "@"
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.916'
$ws.Range("E42").Value = '  -2.18%  '

$rng = $ws.Range("B43:E43")
$rng.NumberFormat = "@"
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '0.9981'
$ws.Range("E43").Value = '  -0.23%  '

$rng = $ws.Range("B44:E44")
$rng.NumberFormat = "@"
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4161'
$ws.Range("E44").Value = '  -1.31%  '

$rng = $ws.Range("B45:E45")
$rng.NumberFormat = "@"
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '101.47'
$ws.Range("E45").Value = '  -4.17%  '

$rng = $ws.Range("B46:E46")
$rng.NumberFormat = "@"
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.355'
$ws.Range("E46").Value = '  -1.40%  '

$rng = $ws.Range("B47:E47")
$rng.NumberFormat = "@"
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1249'
$ws.Range("E47").Value = '  -0.71%  '

$rng = $ws.Range("B48:E48")
$rng.NumberFormat = "@"
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05643'
$ws.Range("E48").Value = '  -0.41%  '

$rng = $ws.Range("B49:E49")
$rng.NumberFormat = "@"
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '32.32'
$ws.Range("E49").Value = '  -0.83%  '

$rng = $ws.Range("B50:E50")
$rng.NumberFormat = "@"
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '55.50'
$ws.Range("E50").Value = '  +0.83%  '

$rng = $ws.Range("B51:E51")
$rng.NumberFormat = "@"
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3677'
$ws.Range("E51").Value = '  -1.25%  '
